$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1426.1428
$ws.Range("J112").Value = 1542.3636
$ws.Range("L112").Value = 4627.0908
$ws.Range("N112").Value = -6843.0908
$ws.Range("H129").Value = 838
$ws.Range("I129").Value = 292.6
$ws.Range("J129").Value = 1097.7142
$ws.Range("K129").Value = 877.8000000000001
$ws.Range("L129").Value = 3293.1426
$ws.Range("M129").Value = 4122.2
$ws.Range("N129").Value = -13293.1426
$ws.Range("H137").Value = 34484756
$ws.Range("I137").Value = 1691.3889
$ws.Range("J137").Value = 90911580
$ws.Range("K137").Value = 5074.1667
$ws.Range("L137").Value = 272734740
$ws.Range("M137").Value = -2524.1667
$ws.Range("N137").Value = -272739840

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10518.573
$ws.Range("I32").Value = 3437.5967
$ws.Range("J32").Value = 26778.592
$ws.Range("K32").Value = 3437.5967
$ws.Range("L32").Value = 26778.592
$ws.Range("M32").Value = -3150.5967
$ws.Range("N32").Value = -27352.592
$ws.Range("H45").Value = 1864.125
$ws.Range("I45").Value = 1848.2858
$ws.Range("K45").Value = 1848.2858
$ws.Range("M45").Value = -1471.2858
$ws.Range("H61").Value = 977.0732
$ws.Range("I61").Value = 887.8333
$ws.Range("J61").Value = 1619.6
$ws.Range("K61").Value = 887.8333
$ws.Range("L61").Value = 1619.6
$ws.Range("M61").Value = -675.8333
$ws.Range("N61").Value = -2043.6
$ws.Range("H74").Value = 1862.2742
$ws.Range("I74").Value = 1486.8085
$ws.Range("J74").Value = 3038.7334
$ws.Range("K74").Value = 1486.8085
$ws.Range("L74").Value = 3038.7334
$ws.Range("M74").Value = -612.8085000000001
$ws.Range("N74").Value = -4786.7334
$ws.Range("H77").Value = 1862.2742
$ws.Range("I77").Value = 1486.8085
$ws.Range("J77").Value = 3038.7334
$ws.Range("K77").Value = 7434.0425
$ws.Range("L77").Value = 15193.667
$ws.Range("M77").Value = -3066.0425
$ws.Range("N77").Value = -23929.667
$ws.Range("H88").Value = 166669570
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 166669570
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 166669570
$ws.Range("N88").Value = -166670382
$ws.Range("H91").Value = 166669570
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 166669570
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 166669570
$ws.Range("N91").Value = -166672378
$ws.Range("H132").Value = 1556.9667
$ws.Range("I132").Value = 1305.2972
$ws.Range("J132").Value = 1961.826
$ws.Range("K132").Value = 3915.8916
$ws.Range("L132").Value = 5885.478
$ws.Range("M132").Value = -1385.8916
$ws.Range("N132").Value = -10945.478
$ws.Range("H136").Value = 977.0732
$ws.Range("I136").Value = 887.8333
$ws.Range("J136").Value = 1619.6
$ws.Range("K136").Value = 2663.4999
$ws.Range("L136").Value = 4858.799999999999
$ws.Range("M136").Value = -113.4998999999998
$ws.Range("N136").Value = -9958.799999999999
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3568.1025
$ws.Range("I20").Value = 1920.3462
$ws.Range("J20").Value = 6863.615
$ws.Range("K20").Value = 1920.3462
$ws.Range("L20").Value = 6863.615
$ws.Range("M20").Value = -1673.3462
$ws.Range("N20").Value = -7357.615
$ws.Range("H86").Value = 18183068
$ws.Range("I86").Value = 18183068
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 18183068
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -18181945
$ws.Range("H89").Value = 18183068
$ws.Range("I89").Value = 18183068
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 90915340
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -90909724
$ws.Range("H134").Value = 1145.8868
$ws.Range("I134").Value = 970.26086
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 2910.78258
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -375.7825800000001
$ws.Range("N134").Value = -11970
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4278504.5
$ws.Range("I31").Value = 6485009
$ws.Range("J31").Value = 3403
$ws.Range("K31").Value = 6485009
$ws.Range("L31").Value = 3403
$ws.Range("M31").Value = -6484714
$ws.Range("N31").Value = -3993
$ws.Range("H34").Value = 4278504.5
$ws.Range("I34").Value = 6485009
$ws.Range("J34").Value = 3403
$ws.Range("K34").Value = 6485009
$ws.Range("L34").Value = 3403
$ws.Range("M34").Value = -6484807
$ws.Range("N34").Value = -3807
$ws.Range("H58").Value = 1436.125
$ws.Range("I58").Value = 815.3570999999999
$ws.Range("J58").Value = 1918.9445
$ws.Range("K58").Value = 815.3570999999999
$ws.Range("L58").Value = 1918.9445
$ws.Range("M58").Value = -612.3570999999999
$ws.Range("N58").Value = -2324.9445
$ws.Range("H122").Value = 6061485.5
$ws.Range("I122").Value = 6061485.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18184456.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -18182006.5
$ws.Range("H132").Value = 1102.3077
$ws.Range("I132").Value = 973
$ws.Range("J132").Value = 1533.3334
$ws.Range("K132").Value = 2919
$ws.Range("L132").Value = 4600.0002
$ws.Range("M132").Value = -389
$ws.Range("N132").Value = -9660.0002
$ws.Range("H134").Value = 2399.4443
$ws.Range("I134").Value = 2495.2144
$ws.Range("J134").Value = 2064.25
$ws.Range("K134").Value = 7485.6432
$ws.Range("L134").Value = 6192.75
$ws.Range("M134").Value = -4950.6432
$ws.Range("N134").Value = -11262.75
$ws.Range("H136").Value = 1436.125
$ws.Range("I136").Value = 815.3570999999999
$ws.Range("J136").Value = 1918.9445
$ws.Range("K136").Value = 2446.0713
$ws.Range("L136").Value = 5756.833500000001
$ws.Range("M136").Value = 103.9287000000004
$ws.Range("N136").Value = -10856.8335
$ws.Range("N122").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 5178.643
$ws.Range("I123").Value = 2197.7778
$ws.Range("J123").Value = 6590.6313
$ws.Range("K123").Value = 6593.3334
$ws.Range("L123").Value = 19771.8939
$ws.Range("M123").Value = -4143.3334
$ws.Range("N123").Value = -24671.8939
$ws.Range("H129").Value = 1279.5
$ws.Range("I129").Value = 1321.7273
$ws.Range("J129").Value = 1124.6666
$ws.Range("K129").Value = 3965.1819
$ws.Range("L129").Value = 3373.9998
$ws.Range("M129").Value = 1034.8181
$ws.Range("N129").Value = -13373.9998
$ws.Range("H130").Value = 1050
$ws.Range("I130").Value = 500
$ws.Range("K130").Value = 1500
$ws.Range("M130").Value = 3520
$ws.Range("H131").Value = 1039.6111
$ws.Range("I131").Value = 525.44446
$ws.Range("J131").Value = 1142.4445
$ws.Range("K131").Value = 1576.33338
$ws.Range("L131").Value = 3427.3335
$ws.Range("M131").Value = 3463.66662
$ws.Range("N131").Value = -13507.3335
$ws.Range("H133").Value = 959.2857
$ws.Range("I133").Value = 505
$ws.Range("K133").Value = 1515
$ws.Range("M133").Value = 3545
$ws.Range("H134").Value = 2205
$ws.Range("I134").Value = 3110
$ws.Range("K134").Value = 9330
$ws.Range("M134").Value = -4260
$ws.Range("H136").Value = 728.75
$ws.Range("I136").Value = 534.7059
$ws.Range("K136").Value = 1604.1177
$ws.Range("M136").Value = 3495.8823
$ws.Range("H137").Value = 3001.76
$ws.Range("I137").Value = 785.8333
$ws.Range("K137").Value = 2357.4999
$ws.Range("M137").Value = 2742.5001
$ws.Range("H138").Value = 791.94446
$ws.Range("I138").Value = 705.5
$ws.Range("K138").Value = 2116.5
$ws.Range("M138").Value = 3023.5
$ws.Range("H139").Value = 4476.2256
$ws.Range("I139").Value = 1743.3077
$ws.Range("J139").Value = 6450
$ws.Range("K139").Value = 5229.9231
$ws.Range("L139").Value = 19350
$ws.Range("M139").Value = -89.92309999999998
$ws.Range("N139").Value = -29630
$ws.Range("H140").Value = 1341.5385
$ws.Range("I140").Value = 1087.7778
$ws.Range("J140").Value = 1912.5
$ws.Range("K140").Value = 3263.3334
$ws.Range("L140").Value = 5737.5
$ws.Range("M140").Value = 1916.6666
$ws.Range("N140").Value = -16097.5
$ws.Range("H141").Value = 1402.1428
$ws.Range("I141").Value = 1107.5
$ws.Range("J141").Value = 1520
$ws.Range("K141").Value = 3322.5
$ws.Range("L141").Value = 4560
$ws.Range("M141").Value = 1857.5
$ws.Range("N141").Value = -14920

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2277.353
$ws.Range("I132").Value = 2150.3333
$ws.Range("J132").Value = 2582.2
$ws.Range("K132").Value = 6450.999899999999
$ws.Range("L132").Value = 7746.599999999999
$ws.Range("M132").Value = -3920.999899999999
$ws.Range("N132").Value = -12806.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 584.1429000000001
$ws.Range("I16").Value = 584.1429000000001
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 584.1429000000001
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -414.1429000000001
$ws.Range("H122").Value = 1864.3
$ws.Range("I122").Value = 1967.0769
$ws.Range("J122").Value = 1196.25
$ws.Range("K122").Value = 5901.2307
$ws.Range("L122").Value = 3588.75
$ws.Range("M122").Value = -3451.2307
$ws.Range("N122").Value = -8488.75
$ws.Range("H132").Value = 9770427
$ws.Range("I132").Value = 20168808
$ws.Range("J132").Value = 2250.697
$ws.Range("K132").Value = 60506424
$ws.Range("L132").Value = 6752.091
$ws.Range("M132").Value = -60503894
$ws.Range("N132").Value = -11812.091
$ws.Range("H136").Value = 3394.2036
$ws.Range("I136").Value = 4606.273
$ws.Range("J136").Value = 1489.5238
$ws.Range("K136").Value = 13818.819
$ws.Range("L136").Value = 4468.5714
$ws.Range("M136").Value = -11268.819
$ws.Range("N136").Value = -9568.571400000001
$ws.Range("N16").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1120.3922
$ws.Range("I132").Value = 808.75
$ws.Range("J132").Value = 2253.6365
$ws.Range("K132").Value = 2426.25
$ws.Range("L132").Value = 6760.9095
$ws.Range("M132").Value = 103.75
$ws.Range("N132").Value = -11820.9095
$ws.Range("H136").Value = 1001.87933
$ws.Range("I136").Value = 518.0513
$ws.Range("K136").Value = 1554.1539
$ws.Range("M136").Value = 995.8461000000002
